$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados..." timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 00:38"

# --- Rank swaps: country overtakes its neighbour in the table (sorted by total cases) ---
# Guinea-Bisau (row136) <-> Angola (row137): Angola grows and moves up to row136
$ws.Range("A136").Value = "Angola"
$ws.Range("A137").Value = "Guinea-Bisau"

# Botsuana (row152) <-> Burkina Faso (row153): Burkina Faso grows and moves up to row152
$ws.Range("A152").Value = "Burkina Faso"
$ws.Range("A153").Value = "Botsuana"

# Republica del Chad (row161) <-> Trinidad yTobago (row162): Trinidad moves up to row161
$ws.Range("A161").Value = "Trinidad yTobago"
$ws.Range("A162").Value = "Republica del Chad"

# Timor Oriental (row202) <-> Santa Lucia (row203): Santa Lucia moves up to row202 (tied totals)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Updated case/death statistics for affected rows ---
$ws.Range("B4").Value = 5872196
$ws.Range("C4").Value = 30768
$ws.Range("D4").Value = 3161135
$ws.Range("E4").Value = 2530498
$ws.Range("G4").Value = 389
$ws.Range("H4").Value = 180563
$ws.Range("B9").Value = 594326
$ws.Range("C9").Value = 9090
$ws.Range("D9").Value = 399357
$ws.Range("E9").Value = 167306
$ws.Range("G9").Value = 210
$ws.Range("H9").Value = 27663
$ws.Range("B11").Value = 541147
$ws.Range("C11").Value = 8044
$ws.Range("D11").Value = 374030
$ws.Range("E11").Value = 149801
$ws.Range("G11").Value = 348
$ws.Range("H11").Value = 17316
$ws.Range("B23").Value = 234478
$ws.Range("C23").Value = 621
$ws.Range("D23").Value = 209600
$ws.Range("E23").Value = 15546
$ws.Range("B34").Value = 97340
$ws.Range("C34").Value = 103
$ws.Range("D34").Value = 65927
$ws.Range("E34").Value = 26151
$ws.Range("G34").Value = 19
$ws.Range("H34").Value = 5262
$ws.Range("B44").Value = 68188
$ws.Range("C44").Value = 332
$ws.Range("D44").Value = 56778
$ws.Range("E44").Value = 8816
$ws.Range("G44").Value = 14
$ws.Range("H44").Value = 2594
$ws.Range("B48").Value = 61747
$ws.Range("C48").Value = 1014
$ws.Range("D48").Value = 48550
$ws.Range("E48").Value = 12021
$ws.Range("G48").Value = 7
$ws.Range("H48").Value = 1176
$ws.Range("B53").Value = 52227
$ws.Range("C53").Value = 322
$ws.Range("D53").Value = 38945
$ws.Range("E53").Value = 12280
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 1002
$ws.Range("B54").Value = 49330
$ws.Range("C54").Value = 292
$ws.Range("D54").Value = 45981
$ws.Range("E54").Value = 3165
$ws.Range("D60").Value = 34400
$ws.Range("E60").Value = 3502
$ws.Range("B61").Value = 39564
$ws.Range("C61").Value = 607
$ws.Range("D61").Value = 29966
$ws.Range("E61").Value = 9269
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 329
$ws.Range("B81").Value = 15287
$ws.Range("C81").Value = 60
$ws.Range("D81").Value = 10338
$ws.Range("E81").Value = 4404
$ws.Range("B90").Value = 10323
$ws.Range("C90").Value = 24
$ws.Range("E90").Value = 909
$ws.Range("B119").Value = 3509
$ws.Range("C119").Value = 54
$ws.Range("D119").Value = 2540
$ws.Range("E119").Value = 932
$ws.Range("B136").Value = 2171
$ws.Range("C136").Value = 37
$ws.Range("D136").Value = 818
$ws.Range("E136").Value = 1257
$ws.Range("G136").Value = 2
$ws.Range("H136").Value = 96
$ws.Range("B137").Value = 2149
$ws.Range("D137").Value = 1015
$ws.Range("E137").Value = 1101
$ws.Range("H137").Value = 33
$ws.Range("B152").Value = 1320
$ws.Range("C152").Value = 23
$ws.Range("D152").Value = 1043
$ws.Range("E152").Value = 222
$ws.Range("H152").Value = 55
$ws.Range("B153").Value = 1308
$ws.Range("D153").Value = 136
$ws.Range("E153").Value = 1169
$ws.Range("H153").Value = 3
$ws.Range("B155").Value = 1277
$ws.Range("C155").Value = 2
$ws.Range("D155").Value = 910
$ws.Range("E155").Value = 340
$ws.Range("B161").Value = 1007
$ws.Range("C161").Value = 77
$ws.Range("D161").Value = 165
$ws.Range("E161").Value = 828
$ws.Range("G161").Value = 1
$ws.Range("H161").Value = 14
$ws.Range("B162").Value = 986
$ws.Range("C162").Value = 4
$ws.Range("D162").Value = 870
$ws.Range("E162").Value = 40
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 76
